$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Tema_1 (homework 1) score for student in row 123 (C123): 0 -> 21 raw score
$ws.Range("C123").Formula = "=MIN(ROUNDUP(21*1.1,0),78)"
# Re-fit the row height so editing the cell doesn't stamp a stray explicit row height
$ws.Rows.Item(123).EntireRow.AutoFit()

# Re-filter the "Grupa" column (column B, the 2nd column of the A1:L134 range)
# from group 241 to group 244 - this drives which rows are shown/hidden
$rng = $ws.Range("A1:L134")
$rng.AutoFilter(2, @("244"), 7)

# Move the active selection to C124, matching where the user ended up working
$ws.Range("C124").Select()
